$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column E
$ws.Range("E1").Value = "strength (RMS)"

# Updated data values for rows 2-19 (columns B, D, E)
$data = @{
    2  = @{ B = 16.4;              D = 14.4;              E = 50.4 }
    3  = @{ B = 15.6;              D = 14.8;              E = 44.4 }
    4  = @{ B = 15.6;              D = 15.2;              E = 40.4 }
    5  = @{ B = 16.8;              D = 14.8;              E = 44.4 }
    6  = @{ B = 18;                D = 12.4;              E = 42 }
    7  = @{ B = 15.6;              D = 15.6;              E = 41.4 }
    8  = @{ B = 16;                D = 14;                E = 41.67 }
    9  = @{ B = 15.2;              D = 15.2;              E = 42.8 }
    10 = @{ B = 16.8;              D = 14;                E = 40.4 }
    11 = @{ B = 16;                D = 14.4;              E = 39.6 }
    12 = @{ B = 16.4;              D = 14;                E = 36.6 }
    13 = @{ B = 16;                D = 14;                E = 44.8 }
    14 = @{ B = 15.6;              D = 15.2;              E = 42.2 }
    15 = @{ B = 15.6;              D = 14.8;              E = 46.4 }
    16 = @{ B = 21.6;              D = 8.800000000000001; E = 35 }
    17 = @{ B = 19.2;              D = 12.8;              E = 43.4 }
    18 = @{ B = 15.5;              D = 15;                E = 48 }
    19 = @{ B = 15.6;              D = 14.8;              E = 47.4 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
}
